$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 2313.611
$ws.Range("M40").Value = -2628.5
$ws.Range("K40").Value = 2803.5
$ws.Range("I40").Value = 2803.5
$ws.Range("N113").Value = -9878.3333
$ws.Range("H113").Value = 3402.75
$ws.Range("L113").Value = 3370.3333
$ws.Range("J113").Value = 3370.3333
$ws.Range("K116").Value = 3499.5
$ws.Range("L116").Value = 0
$ws.Range("J116").Value = 0
$ws.Range("N116").Value = ""
$ws.Range("M116").Value = -57.5
$ws.Range("I116").Value = 3499.5
$ws.Range("H116").Value = 3499.5
$ws.Range("M132").Value = -1635.5681
$ws.Range("I132").Value = 1388.5227
$ws.Range("H132").Value = 1591.9
$ws.Range("K132").Value = 4165.5681
$ws.Range("N138").Value = -20729.738
$ws.Range("H138").Value = 2606.21
$ws.Range("I138").Value = 977.4286
$ws.Range("K138").Value = 2932.2858
$ws.Range("J138").Value = 3483.246
$ws.Range("M138").Value = 2207.7142
$ws.Range("L138").Value = 10449.738
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("K86").Value = 0
$ws.Range("J86").Value = 40314
$ws.Range("I86").Value = 0
$ws.Range("M86").Value = ""
$ws.Range("L86").Value = 40314
$ws.Range("H86").Value = 40314
$ws.Range("N86").Value = -42686
$ws.Range("L89").Value = 120942
$ws.Range("J89").Value = 40314
$ws.Range("K89").Value = 0
$ws.Range("I89").Value = 0
$ws.Range("M89").Value = ""
$ws.Range("H89").Value = 40314
$ws.Range("N89").Value = -132798
$ws.Range("H105").Value = 0
$ws.Range("L105").Value = 0
$ws.Range("J105").Value = 0
$ws.Range("N105").Value = ""
$ws.Range("I122").Value = 1944.5
$ws.Range("M122").Value = -3383.5
$ws.Range("H122").Value = 12502067
$ws.Range("K122").Value = 5833.5
$ws.Range("N124").Value = ""
$ws.Range("J124").Value = 0
$ws.Range("H124").Value = 0
$ws.Range("L124").Value = 0
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("N58").Value = -50388
$ws.Range("H58").Value = 49800
$ws.Range("L58").Value = 49800
$ws.Range("J58").Value = 49800
$ws.Range("J59").Value = 59775
$ws.Range("M59").Value = ""
$ws.Range("L59").Value = 59775
$ws.Range("H59").Value = 59775
$ws.Range("I59").Value = 0
$ws.Range("K59").Value = 0
$ws.Range("N59").Value = -61469
$ws.Range("N138").Value = -44113
$ws.Range("H138").Value = 33833
$ws.Range("J138").Value = 33833
$ws.Range("L138").Value = 33833
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("J48").Value = 12336.714
$ws.Range("H48").Value = 12336.714
$ws.Range("L48").Value = 12336.714
$ws.Range("N48").Value = -13288.714
$ws.Range("M58").Value = -7578560.5
$ws.Range("N58").Value = -4425.625
$ws.Range("I58").Value = 7578763.5
$ws.Range("H58").Value = 3250338.5
$ws.Range("L58").Value = 4019.625
$ws.Range("J58").Value = 4019.625
$ws.Range("K58").Value = 7578763.5
$ws.Range("J120").Value = 0
$ws.Range("N120").Value = ""
$ws.Range("H120").Value = 0
$ws.Range("L120").Value = 0
$ws.Range("H121").Value = 24000
$ws.Range("N121").Value = -26620
$ws.Range("L121").Value = 24000
$ws.Range("J121").Value = 24000
$ws.Range("I122").Value = 3924.5334
$ws.Range("M122").Value = -9323.600199999999
$ws.Range("H122").Value = 8667.179
$ws.Range("J122").Value = 14139.462
$ws.Range("L122").Value = 42418.386
$ws.Range("N122").Value = -47318.386
$ws.Range("K122").Value = 11773.6002
$ws.Range("M123").Value = -35100
$ws.Range("J123").Value = 53918
$ws.Range("N123").Value = -63718
$ws.Range("L123").Value = 53918
$ws.Range("H123").Value = 52652.727
$ws.Range("I123").Value = 40000
$ws.Range("K123").Value = 40000
$ws.Range("I136").Value = 7578763.5
$ws.Range("K136").Value = 22736290.5
$ws.Range("H136").Value = 3250338.5
$ws.Range("J136").Value = 4019.625
$ws.Range("L136").Value = 12058.875
$ws.Range("N136").Value = -17158.875
$ws.Range("M136").Value = -22733740.5
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("K4").Value = 1800
$ws.Range("M4").Value = -1688
$ws.Range("H4").Value = 1081.8182
$ws.Range("N4").Value = -5204
$ws.Range("J4").Value = 1660
$ws.Range("L4").Value = 4980
$ws.Range("I4").Value = 600
$ws.Range("M68").Value = -1229
$ws.Range("L68").Value = 27899.25
$ws.Range("I68").Value = 680
$ws.Range("H68").Value = 6764.5293
$ws.Range("K68").Value = 2040
$ws.Range("N68").Value = -29521.25
$ws.Range("J68").Value = 9299.75
$ws.Range("L71").Value = 83697.75
$ws.Range("I71").Value = 680
$ws.Range("N71").Value = -91809.75
$ws.Range("M71").Value = -2064
$ws.Range("H71").Value = 6764.5293
$ws.Range("K71").Value = 6120
$ws.Range("J71").Value = 9299.75
$ws.Range("H105").Value = 6765
$ws.Range("L105").Value = 20295
$ws.Range("J105").Value = 6765
$ws.Range("N105").Value = -25537
$ws.Range("J109").Value = 2914.2104
$ws.Range("K109").Value = 641.00001
$ws.Range("M109").Value = 398.99999
$ws.Range("I109").Value = 213.66667
$ws.Range("H109").Value = 2266.08
$ws.Range("L109").Value = 8742.6312
$ws.Range("N109").Value = -10822.6312
$ws.Range("N113").Value = -6519.7142
$ws.Range("M113").Value = -156.1802899999998
$ws.Range("H113").Value = 762.89026
$ws.Range("L113").Value = 2179.7142
$ws.Range("I113").Value = 775.39343
$ws.Range("K113").Value = 2326.18029
$ws.Range("J113").Value = 726.5714
$ws.Range("L131").Value = 2545.857
$ws.Range("N131").Value = -12625.857
$ws.Range("H131").Value = 527.87
$ws.Range("K131").Value = 886.81035
$ws.Range("I131").Value = 295.60345
$ws.Range("J131").Value = 848.619
$ws.Range("M131").Value = 4153.18965
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H15").Value = 23000
$ws.Range("J15").Value = 23000
$ws.Range("N15").Value = -23576
$ws.Range("L15").Value = 23000
$ws.Range("N81").Value = -24996
$ws.Range("J81").Value = 23000
$ws.Range("H81").Value = 23000
$ws.Range("L81").Value = 23000
$ws.Range("L84").Value = 69000
$ws.Range("J84").Value = 23000
$ws.Range("H84").Value = 23000
$ws.Range("N84").Value = -78984
$ws.Range("H88").Value = 31597.5
$ws.Range("N88").Value = -32499.5
$ws.Range("J88").Value = 31597.5
$ws.Range("L88").Value = 31597.5
$ws.Range("H91").Value = 31597.5
$ws.Range("L91").Value = 31597.5
$ws.Range("J91").Value = 31597.5
$ws.Range("N91").Value = -34717.5
$ws.Range("I122").Value = 5492.5835
$ws.Range("M122").Value = -14027.7505
$ws.Range("H122").Value = 4987.4
$ws.Range("J122").Value = 2966.6667
$ws.Range("L122").Value = 8900.000100000001
$ws.Range("N122").Value = -13800.0001
$ws.Range("K122").Value = 16477.7505
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("L2").Value = 82000
$ws.Range("N2").Value = -82224
$ws.Range("H2").Value = 82000
$ws.Range("J2").Value = 82000
$ws.Range("I122").Value = 5135.2383
$ws.Range("M122").Value = -12955.7149
$ws.Range("H122").Value = 6688
$ws.Range("J122").Value = 10311.111
$ws.Range("L122").Value = 30933.333
$ws.Range("N122").Value = -35833.333
$ws.Range("K122").Value = 15405.7149
$ws.Range("L132").Value = 11233.7499
$ws.Range("M132").Value = -2588.1305
$ws.Range("I132").Value = 1706.0435
$ws.Range("N132").Value = -16293.7499
$ws.Range("H132").Value = 2404.9714
$ws.Range("K132").Value = 5118.1305
$ws.Range("J132").Value = 3744.5833
$ws.Range("L135").Value = 57386.445
$ws.Range("H135").Value = 57386.445
$ws.Range("J135").Value = 57386.445
$ws.Range("N135").Value = -67526.44500000001
$ws.Range("I136").Value = 2385.3125
$ws.Range("K136").Value = 7155.9375
$ws.Range("H136").Value = 3705.3389
$ws.Range("J136").Value = 9465.454
$ws.Range("L136").Value = 28396.362
$ws.Range("N136").Value = -33496.362
$ws.Range("M136").Value = -4605.9375
$ws.Range("J139").Value = 59150
$ws.Range("H139").Value = 59150
$ws.Range("L139").Value = 59150
$ws.Range("N139").Value = -69430
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("L135").Value = 222249380
$ws.Range("H135").Value = 222249380
$ws.Range("J135").Value = 222249380
$ws.Range("N135").Value = -222259520
$ws.Range("J140").Value = 51595
$ws.Range("N140").Value = -61955
$ws.Range("H140").Value = 51595
$ws.Range("L140").Value = 51595
